$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report used to list three hyperlinked stats (Participants/Pages/Modules),
# each linking out to the Canvas course, plus a URL row (B4) that keeps its
# link. "Pages"/"Modules" get renamed, a new "Assessment count" row is
# inserted, and "Time taken to generate" shifts down to row 9 with a new
# value.

# Worksheet.Hyperlinks.Delete() removes every hyperlink on the sheet at
# once (it isn't scoped to the calling Range), so drop all four here and
# re-create only the one that should survive (B4).
$ws.Range("B4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B4"), "https://wisdomlearning.instructure.com/courses/704") | Out-Null
$ws.Range("B4").Style = "Hyperlink"

# --- Row 5: "Participants" - plain text "0" (no more hyperlink styling) ---
$ws.Range("B5").Style = "Normal"
$ws.Range("B5").Value = "'0"

# --- Row 6: "Pages" -> "Page count", value becomes numeric 30 ---
$ws.Range("A6").Value = "Page count"
$ws.Range("B6").Style = "Normal"
$ws.Range("B6").Value = 30

# --- Row 7: "Modules" -> "Module count", value becomes numeric 6 ---
$ws.Range("A7").Value = "Module count"
$ws.Range("B7").Style = "Normal"
$ws.Range("B7").Value = 6

# --- New row 8: "Assessment count" = 16 ---
$ws.Range("A8").Value = "Assessment count"
$ws.Range("B8").Value = 16

# --- Row 9 (previously row 7's "Time taken to generate"): shift down, new value ---
$ws.Range("A9").Value = "Time taken to generate"

# B9 holds numeric-looking text ("74.2") that must stay text rather than
# become a real number, so enter it with a leading apostrophe.
$ws.Range("B9").Value = "'74.2"

# B5/B9 were entered with a quote-prefix just to force text storage; drop
# that quote-prefix formatting now that the value is stored, returning both
# cells to the plain (unstyled) Normal look the report uses elsewhere.
$ws.Range("B5").Style = "Normal"
$ws.Range("B9").Style = "Normal"

Write-Output "edit complete"
